$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the card's text rows (A2:A7) into a single Python-tuple-repr string in A2
$ws.Range("A2").Value = '(''Sword of Kaldra'', [''{4}'', ''Legendary Artifact — Equipment'', ''Equipped creature gets +5/+5.'', ''Whenever equipped creature deals damage to a creature, exile that creature. (Exile it only if it’s still on the battlefield.)'', ''Equip {4} ({4}: Attach to target creature you control. Equip only as a sorcery. This card enters the battlefield unattached and stays on the battlefield if the creature leaves.)''])'

# Clear the now-merged rows 3-7 so the sheet's used range shrinks to A1:A2
$ws.Range("A3:A7").ClearContents()
